$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45986, 11416.5763447824, 10266.4528647254, 17427.4, 6644.91131787776, -21.5014923915335),
    @(45987, 11320.760651225, 10549.1791318776, 9003.4, 6631.19271108869, 340.707160123595),
    @(45988, 11415.4915409204, 10861.5542150527, 9003.4, 6731.03079451093, 357.882708731819),
    @(45989, 11326.7043167823, 10317.4630871922, 9003.4, 6686.16489975834, 333.342832789604),
    @(45990, 4509.24399740186, 7357.30763749824, 9003.4, 6144.83142324198, 187.447460864176),
    @(45991, 4285.59691763933, 7038.10284369069, 9003.4, 6023.17493479135, 169.078240770085),
    @(45992, 11775.3234136769, 10900.8489534046, 11064.26, 7488.48046616165, 305.211225815259),
    @(45993, 11775.3234136769, 11119.6314722592, 11064.26, 7488.48046616165, 314.327164100867),
    @(45994, 11775.3234136769, 10874.060588647, 11064.26, 7488.48046616165, 304.09504395036),
    @(45995, 11775.3234136769, 10852.3324063833, 11064.26, 7488.48046616165, 303.189703022705),
    @(45996, 11775.3234136769, 10096.3151674463, 11064.26, 7488.48046616165, 271.688984733663),
    @(45997, 4605.75571869631, 6806.13023818529, 11064.26, 7098.30328362905, 118.340563408931),
    @(45998, 4450.60726673227, 6530.13971821467, 11064.26, 7083.23595700977, 106.213153134352),
    @(45999, 4752.2157406543, 6651.40742507455, 11064.26, 7516.95884783414, 129.337761371195)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
